# Update with token state in DB — open-source the generator codes.
#
# 1. Add a new "token" worksheet (uuid / username / role), based on the
#    "acl" sheet so the dimension (A1:H7) and the two list data
#    validations (B5:B7, E5:H7) come along for free.
# 2. Reorder tabs: config, account, token, acl, audit, test.
# 3. account.password length 16 -> 128.
# 4. Misc selection / active-tab bookkeeping to match the authored state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the "token" sheet from a copy of "acl" (same 3 data-row shape,
#    same two list data validations already scoped to rows 5-7)
# ---------------------------------------------------------------------
$acl = $wb.Worksheets.Item("acl")
[void]$acl.Copy($null, $acl)
$token = $wb.Worksheets.Item("acl (2)")
$token.Name = "token"

# Row 5: uuid / VARCHAR / 64 / Key=Yes, Index=No, Nullable=No, Unsigned=No
# (written before B2 so the two brand-new shared strings are interned in
# "uuid", "token" order, same as the authored file)
$token.Range("A5").Value = "uuid"
$token.Range("C5").Value = 64
[void]$token.Range("D5").ClearContents()
$token.Range("E5").Value = "Yes"
$token.Range("F5").Value = "No"
$token.Range("G5").Value = "No"
$token.Range("H5").Value = "No"

# Table-name cell (B2)
$token.Range("B2").Value = "token"

# Row 6: username / VARCHAR / 32 / all No
$token.Range("A6").Value = "username"
$token.Range("C6").Value = 32
[void]$token.Range("D6").ClearContents()
$token.Range("E6").Value = "No"
$token.Range("F6").Value = "No"
$token.Range("G6").Value = "No"
$token.Range("H6").Value = "No"

# Row 7: role / VARCHAR / 16 / all No
$token.Range("A7").Value = "role"
$token.Range("C7").Value = 16
[void]$token.Range("D7").ClearContents()
$token.Range("E7").Value = "No"
$token.Range("F7").Value = "No"
$token.Range("G7").Value = "No"
$token.Range("H7").Value = "No"

# ---------------------------------------------------------------------
# 2. Reorder the tabs: config, account, token, acl, audit, test
# ---------------------------------------------------------------------
$config = $wb.Worksheets.Item("config")
$account = $wb.Worksheets.Item("account")
[void]$account.Move($null, $config)

# token currently sits right after acl; move it to sit right after account
$account = $wb.Worksheets.Item("account")
$token = $wb.Worksheets.Item("token")
[void]$token.Move($null, $account)

# ---------------------------------------------------------------------
# 3. account.password length 16 -> 128
# ---------------------------------------------------------------------
$account = $wb.Worksheets.Item("account")
$account.Range("C6").Value = 128

# ---------------------------------------------------------------------
# 4. Per-sheet selections (each Range.Select() also activates that
#    sheet, so these run in an order that leaves "account" active last,
#    matching the authored activeTab = 1 / tabSelected on "account")
# ---------------------------------------------------------------------
$config = $wb.Worksheets.Item("config")
[void]$config.Range("G8").Select()

$token = $wb.Worksheets.Item("token")
[void]$token.Range("B2").Select()

$acl = $wb.Worksheets.Item("acl")
[void]$acl.Range("L14").Select()

# audit is untouched (selection stays B2, it was never the active sheet)

$test = $wb.Worksheets.Item("test")
[void]$test.Range("F16").Select()

# account selected last -> becomes (and stays) the active tab
$account = $wb.Worksheets.Item("account")
[void]$account.Range("C6").Select()
